$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, 2, 1),
    @(4, 4, 1),
    @(2, 1, 2),
    @(1, 3, 1),
    @(4, 1, 2),
    @(3, 3, 1),
    @(2, 4, 1),
    @(1, 2, 1),
    @(3, 1, 2),
    @(4, 3, 1),
    @(2, 2, 1),
    @(1, 4, 1),
    @(4, 2, 2),
    @(2, 3, 1),
    @(1, 1, 1),
    @(3, 4, 1),
    @(3, 2, 2),
    @(4, 4, 1),
    @(2, 1, 2),
    @(1, 3, 1),
    @(4, 1, 2),
    @(3, 3, 1),
    @(2, 4, 1),
    @(1, 2, 1),
    @(3, 1, 2),
    @(4, 3, 1),
    @(2, 2, 1),
    @(1, 4, 1),
    @(4, 2, 1),
    @(2, 3, 1),
    @(1, 1, 1),
    @(3, 4, 1),
)

$startRow = 24
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$ws.Range("E17").Select()
